$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.263.55'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.85%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.243.99'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.39'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('E6').Value = '  -2.97%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '74.58'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.74%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.617'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.23%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.07'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +6.65%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0939'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.30%  '
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('E13').Value = '  -2.53%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.55'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.14%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.851'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.231.00'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.90%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '42.095.89'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0₃0989'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '72.02'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.24'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '231.89'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.72'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +35.77%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.32'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('E26').Value = '  -3.17%  '
$ws.Range('E27').Value = '  -2.57%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '169.53'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.08'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.78%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.68'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('E31').Value = '  -4.78%  '
$ws.Range('E32').Value = '  -3.14%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '30.52'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('E34').Value = '  -1.41%  '
$ws.Range('E35').Value = '  +9.73%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.51'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '13.67'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '62.26'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.54%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.203'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '106.65'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('E44').Value = '  +1.86%  '
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.997'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.36'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -6.67%  '
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.26'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.10'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.09%  '
